$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (index 1): update "想去人数" (F column) counters
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 123
$ws1.Range("F3").Value = 949
$ws1.Range("F4").Value = 592
$ws1.Range("F5").Value = 2819
$ws1.Range("F6").Value = 771
$ws1.Range("F8").Value = 588
$ws1.Range("F9").Value = 75
$ws1.Range("F12").Value = 423
$ws1.Range("F15").Value = 1246
$ws1.Range("F18").Value = 2654
$ws1.Range("F23").Value = 512
$ws1.Range("F24").Value = 557
$ws1.Range("F32").Value = 4636
$ws1.Range("F33").Value = 211

# ---------------------------------------------------------------------------
# Sheet "演出" (index 2): update "想去人数" (F column) counters
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F5").Value = 4
$ws2.Range("F7").Value = 325
$ws2.Range("F12").Value = 162
$ws2.Range("F21").Value = 270
$ws2.Range("F23").Value = 292
$ws2.Range("F31").Value = 9
$ws2.Range("F32").Value = 496
$ws2.Range("F33").Value = 12

# ---------------------------------------------------------------------------
# Sheet "本地生活" (index 3): the "上海・明日方舟主题店·[SWEET ZONE甜蜜区域]"
# listing (row 5) has expired and fallen off the feed. Rows 6-8 shift up into
# rows 5-7 (column A, the fixed sequence number, is left untouched), and the
# now-empty trailing row 8 is removed.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("F4").Value = 1445

# Row 5 <= old row 6 (NIJISANJI EN)
$ws3.Range("B5").Value = "'2024-03-21"
$ws3.Range("B5").ClearFormats()
$ws3.Range("C5").Value = "上海·NIJISANJI EN 官方授权主题店"
$ws3.Range("D5").Value = "西藏北路166号（地铁8号线曲阜路下） 静安大悦城"
$ws3.Range("E5").Value = "2024.03.21 00:00-04.28 23:59"
$ws3.Range("F5").Value = 564
$ws3.Range("G5").Value = 30
$ws3.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=82858"
$ws3.Range("I5").Value = "//i1.hdslb.com/bfs/openplatform/202403/eeG6Usri1710399065622.jpeg"

# Row 6 <= old row 7 (Paradox Live)
$ws3.Range("B6").Value = "'2024-04-04"
$ws3.Range("B6").ClearFormats()
$ws3.Range("C6").Value = "上海·Paradox Live·[灯光渐强]"
$ws3.Range("D6").Value = "南京东路830号第一百货商业中心B馆5楼(海底捞旁边) 第一百货商业中心"
$ws3.Range("E6").Value = "2024.04.04 00:00-05.20 23:59"
$ws3.Range("F6").Value = 196
$ws3.Range("G6").Value = 30
$ws3.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=83450"
$ws3.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202403/eDK8grIC1711591902400.jpeg"

# Row 7 <= old row 8 (NIJISANJI EN x animate cafe), want-to-go count +1
$ws3.Range("B7").Value = "'2024-04-24"
$ws3.Range("B7").ClearFormats()
$ws3.Range("C7").Value = "上海·「NIJISANJI EN x animate cafe」"
$ws3.Range("D7").Value = "西藏北路198号大悦城北座8楼N809-1 animate cafe上海店"
$ws3.Range("E7").Value = "2024.04.24 00:00-05.22 23:59"
$ws3.Range("F7").Value = 228
$ws3.Range("G7").Value = 30
$ws3.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=83223"
$ws3.Range("I7").Value = "//i0.hdslb.com/bfs/openplatform/202403/LzJJK9lc1711096202393.jpeg"

# The old row 8 content has now been copied into row 7; drop the stale last row.
$ws3.Rows.Item(8).Delete()

# ---------------------------------------------------------------------------
# Sheet "全部类型" (index 4): update "想去人数" (F column) counters
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 1445
$ws4.Range("F5").Value = 123
$ws4.Range("F8").Value = 949
$ws4.Range("F9").Value = 592
$ws4.Range("F10").Value = 2819
$ws4.Range("F11").Value = 771
$ws4.Range("F13").Value = 588
$ws4.Range("F14").Value = 75
$ws4.Range("F17").Value = 423
$ws4.Range("F18").Value = 325
$ws4.Range("F23").Value = 1246
$ws4.Range("F27").Value = 2654
$ws4.Range("F35").Value = 228
$ws4.Range("F36").Value = 512
$ws4.Range("F37").Value = 557
$ws4.Range("F38").Value = 557
$ws4.Range("F39").Value = 270
$ws4.Range("F42").Value = 292
$ws4.Range("F47").Value = 4636
$ws4.Range("F48").Value = 211
$ws4.Range("F50").Value = 496
